$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.015.85"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "1.883.72"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("D4").Value = "'0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'0.7371"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").Value = "'241.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.9999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E8").Value = "  +1.14%  "
$ws.Range("D9").Value = "'0.07165"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.77%  "
$ws.Range("D10").Value = "'24.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.24%  "
$ws.Range("D11").Value = "'0.08313"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.33%  "
$ws.Range("D12").Value = "'0.7556"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.402"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.891.38"
$ws.Range("E14").Value = "  +5.62%  "
$ws.Range("D15").Value = "'92.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("B16").Value = "Uniswap"
$ws.Range("C16").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D16").Value = "'6.148"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "30.053.01"
$ws.Range("E17").Value = "  +0.87%  "
$ws.Range("D18").Value = "'249.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.62%  "
$ws.Range("D19").Value = "'13.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").Value = "'0.000007853"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").Value = "'7.905"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("D23").Value = "'0.9992"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "'0.1569"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.27%  "
$ws.Range("D25").Value = "'9.264"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("D26").Value = "'164.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("D27").Value = "'18.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "'2.046"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").Value = "'4.554"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").Value = "'4.178"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.48%  "
$ws.Range("D33").Value = "'0.05322"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("E34").Value = "  +0.73%  "
$ws.Range("D35").Value = "'0.7681"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.06%  "
$ws.Range("D36").Value = "'0.9999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("D37").Value = "'2.726"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("D38").Value = "'0.01955"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").Value = "'0.4562"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.21%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.8813"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.53%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'6.041"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'72.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.083.78"
$ws.Range("E44").Value = "  -1.51%  "
$ws.Range("D45").Value = "'104.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.86%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").Value = "'7.530"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.54%  "
$ws.Range("D49").Value = "2.064.86"
$ws.Range("E49").Value = "  +1.31%  "
$ws.Range("D50").Value = "'9.533"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.98%  "
$ws.Range("D51").Value = "'2.894"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.55%  "
